$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 4 blank rows right after row 4 (the "회원" header row) so the new
#    member-API rows (join / signup / register / registerSeller) have space.
#    Old rows 5,6,7 (상품/고객센터/관리자) shift down to 9,10,11.
# ---------------------------------------------------------------------------
$ws.Range("A5:A8").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Remove every existing hyperlink so we can re-add them cleanly once all
#    the cell values are in their final place (row-insert does not shift the
#    <hyperlinks> anchors, so re-creating them avoids stale references).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 3. Row 4 (existing "회원" / login row) gets new content.
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "http://localhost:8080/Kmarket/member/login.do"
$ws.Range("C4").Value = "POST"
$ws.Range("D4").Value = "회원 로그인 화면`n- 로그인/세션에 정보 저장`n- 자동 로그인/쿠키 생성`n- 아이디 찾기`n- 비밀번호 찾기/비밀번호 변경`n- 회원가입 페이지로 이동`n- 헤더 링크 연결"
$ws.Range("E4").Value = "조수빈"

# ---------------------------------------------------------------------------
# 4. New rows 5-8: member join / signup / register / registerSeller.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "http://localhost:8080/Kmarket/member/join.do"
$ws.Range("C5").Value = "GET"
$ws.Range("D5").Value = "회원 가입 페이지로 연결하는 링크`n- 개인 구매회원 회원가입 페이지로 연결`n- 판매회원 회원가입 페이지로 연결`n- 헤더 링크 연결"
$ws.Range("E5").Value = "조수빈"

$ws.Range("B6").Value = "http://localhost:8080/Kmarket/member/signup.do"
$ws.Range("C6").Value = "GET"
$ws.Range("D6").Value = "약관 동의 화면`n- 약관 내용 로드`n- 필수 동의 요소 모두 동의 체크했는지 확인하기`n- (개인/판매) 회원 가입 페이지로 연결`n- 헤더 링크 연결"
$ws.Range("E6").Value = "조수빈"

$ws.Range("B7").Value = "http://localhost:8080/Kmarket/member/register.do"
$ws.Range("C7").Value = "POST"
$ws.Range("D7").Value = "개인 회원 회원가입 화면`n- 입력한 정보 유효성 및 중복 검사`n- 필수 정보 모두 기입했는지 검사`n- 우편번호 찾기 페이지 구현`n- 회원가입 완료 후 메인 화면으로 리다이렉트`n- 헤더 링크 구현"
$ws.Range("E7").Value = "조수빈"

$ws.Range("B8").Value = "http://localhost:8080/Kmarket/member/registerSeller.do"
$ws.Range("C8").Value = "POST"
$ws.Range("D8").Value = "판매 회원 회원가입 화면`n- 입력한 정보 유효성 및 중복 검사`n- 필수 정보 모두 기입했는지 검사`n- 우편번호 찾기 페이지 구현`n- 회원가입 완료 후 메인 화면으로 리다이렉트`n- 헤더 링크 구현"
$ws.Range("E8").Value = "조수빈"

# ---------------------------------------------------------------------------
# 5. Merge the "회원" label across the whole member-API block.
# ---------------------------------------------------------------------------
$ws.Range("A4:A8").Merge()

# ---------------------------------------------------------------------------
# 6. Row heights for the new / expanded description rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 115.5
$ws.Rows.Item(5).RowHeight = 66
$ws.Rows.Item(6).RowHeight = 82.5
$ws.Rows.Item(7).RowHeight = 99
$ws.Rows.Item(8).RowHeight = 99

# ---------------------------------------------------------------------------
# 7. Wrap text for the long description column on rows 4-8.
# ---------------------------------------------------------------------------
$ws.Range("D4:D8").WrapText = $true
$ws.Range("D4:D8").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 8. Re-create every hyperlink against its final cell location.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "http://localhost:8080/Kmarket/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "http://localhost:8080/Kmarket/index.do") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "http://localhost:8080/Kmarket/member/login.do") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "http://localhost:8080/Kmarket/product/…") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "http://localhost:8080/Kmarket/cs/…") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "http://localhost:8080/Kmarket/admin/…") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "http://localhost:8080/Kmarket/member/join.do") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "http://localhost:8080/Kmarket/member/register.do") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "http://localhost:8080/Kmarket/member/signup.do") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "http://localhost:8080/Kmarket/member/registerSeller.do") | Out-Null

# ---------------------------------------------------------------------------
# 9. Sheet view: scroll so row 7 is at the top, selection resting on D8.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("D8").Select()
